$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.171.03'
$ws.Range("E2").Value = '  +0.27%  '

# Row 3
$ws.Range("D3").Value = '1.903.07'
$ws.Range("E3").Value = '  +0.65%  '

# Row 4
$ws.Range("E4").Value = '  +0.25%  '

# Row 5
$ws.Range("D5").Value = '''306.29'
$ws.Range("E5").Value = '  -0.45%  '

# Row 6
$ws.Range("E6").Value = '  +0.19%  '

# Row 7
$ws.Range("D7").Value = '''0.5256'
$ws.Range("E7").Value = '  +1.29%  '

# Row 8
$ws.Range("D8").Value = '''0.3776'
$ws.Range("E8").Value = '  +1.34%  '

# Row 9
$ws.Range("D9").Value = '''0.07260'
$ws.Range("E9").Value = '  +0.69%  '

# Row 10
$ws.Range("E10").Value = '  +0.06%  '

# Row 11
$ws.Range("D11").Value = '''0.8992'
$ws.Range("E11").Value = '  -0.65%  '

# Row 12
$ws.Range("D12").Value = '''0.08379'
$ws.Range("E12").Value = '  +9.45%  '

# Row 13
$ws.Range("D13").Value = '1.892.86'
$ws.Range("E13").Value = '  +0.26%  '

# Row 14
$ws.Range("D14").Value = '''94.85'
$ws.Range("E14").Value = '  -0.47%  '

# Row 15
$ws.Range("D15").Value = '''5.271'
$ws.Range("E15").Value = '  -0.13%  '

# Row 16
$ws.Range("E16").Value = '  +0.30%  '

# Row 17
$ws.Range("D17").Value = '''0.000008619'
$ws.Range("E17").Value = '  +1.28%  '

# Row 18
$ws.Range("D18").Value = '''14.55'
$ws.Range("E18").Value = '  +1.60%  '

# Row 19
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  +0.16%  '

# Row 20
$ws.Range("D20").Value = '27.207.87'
$ws.Range("E20").Value = '  +0.26%  '

# Row 21
$ws.Range("D21").Value = '''5.063'
$ws.Range("E21").Value = '  +0.17%  '

# Row 22
$ws.Range("D22").Value = '2.135.79'
$ws.Range("E22").Value = '  +1.52%  '

# Row 23
$ws.Range("E23").Value = '  +0.38%  '

# Row 24
$ws.Range("D24").Value = '''6.434'
$ws.Range("E24").Value = '  -0.48%  '

# Row 25
$ws.Range("D25").Value = '''146.62'
$ws.Range("E25").Value = '  +0.56%  '

# Row 26
$ws.Range("E26").Value = '  +6.71%  '

# Row 27
$ws.Range("E27").Value = '  -1.69%  '

# Row 28
$ws.Range("D28").Value = '''18.14'
$ws.Range("E28").Value = '  +0.47%  '

# Row 29
$ws.Range("D29").Value = '''114.84'
$ws.Range("E29").Value = '  +0.17%  '

# Row 30
$ws.Range("D30").Value = '''4.931'
$ws.Range("E30").Value = '  -0.25%  '

# Row 31
$ws.Range("D31").Value = '''4.794'
$ws.Range("E31").Value = '  -0.09%  '

# Row 32
$ws.Range("D32").Value = '''0.09281'
$ws.Range("E32").Value = '  +0.81%  '

# Row 33
$ws.Range("D33").Value = '''0.8160'
$ws.Range("E33").Value = '  +6.93%  '

# Row 34
$ws.Range("D34").Value = '''0.05056'
$ws.Range("E34").Value = '  +0.08%  '

# Row 35
$ws.Range("D35").Value = '''1.239'
$ws.Range("E35").Value = '  +3.47%  '

# Row 36
$ws.Range("E36").Value = '  -2.25%  '

# Row 37
$ws.Range("E37").Value = '  +2.12%  '

# Row 38
$ws.Range("D38").Value = '''2.590'
$ws.Range("E38").Value = '  +0.91%  '

# Row 39
$ws.Range("D39").Value = '''0.5718'
$ws.Range("E39").Value = '  +1.76%  '

# Row 40
$ws.Range("D40").Value = '''0.01985'
$ws.Range("E40").Value = '  -0.50%  '

# Row 41
$ws.Range("D41").Value = '''1.069'
$ws.Range("E41").Value = '  -0.63%  '

# Row 42
$ws.Range("D42").Value = '''6.671'
$ws.Range("E42").Value = '  +1.10%  '

# Row 43
$ws.Range("D43").Value = '''8.953'
$ws.Range("E43").Value = '  +0.83%  '

# Row 44
$ws.Range("D44").Value = '''117.92'
$ws.Range("E44").Value = '  -0.45%  '

# Row 45
$ws.Range("D45").Value = '''0.1513'
$ws.Range("E45").Value = '  +0.26%  '

# Row 46
$ws.Range("D46").Value = '''0.4838'
$ws.Range("E46").Value = '  +0.72%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''10.20'
$ws.Range("E47").Value = '  +0.51%  '

# Row 48
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '''1.000'
$ws.Range("E48").Value = '  +0.18%  '

# Row 49
$ws.Range("D49").Value = '''1.617'
$ws.Range("E49").Value = '  +2.51%  '

# Row 50
$ws.Range("D50").Value = '''37.47'
$ws.Range("E50").Value = '  +0.89%  '

# Row 51
$ws.Range("D51").Value = '''63.63'
$ws.Range("E51").Value = '  +0.17%  '
